$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50, shifting existing rows 50-79 down to 51-80.
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with its values.
$ws.Range("A50").Value = 10
$ws.Range("B50").Value = "Vega Modelo de Temuco"
$ws.Range("C50").Value = "La Araucanía"
$ws.Range("D50").Value = 45029
$ws.Range("E50").Value = 9
$ws.Range("F50").Value = 100112010
$ws.Range("G50").Value = "Achicoria"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 100
$ws.Range("K50").Value = 10000
$ws.Range("L50").Value = 10000
$ws.Range("M50").Value = 10000
$ws.Range("N50").Value = "$/caja 18 unidades"
$ws.Range("O50").Value = "Región Metropolitana"
$ws.Range("P50").Value = 556
$ws.Range("Q50").Value = 18
$ws.Range("R50").Value = "Hortaliza"
